$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 9
$ws.Range("B2").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C2").Value = 'Metropolitana'
$ws.Range("D2").Value = 44349
$ws.Range("E2").Value = 13
$ws.Range("F2").Value = 100112035
$ws.Range("G2").Value = 'Bruselas (repollito)'
$ws.Range("H2").Value = 'Sin especificar'
$ws.Range("I2").Value = 'Primera'
$ws.Range("J2").Value = 21
$ws.Range("K2").Value = 24000
$ws.Range("L2").Value = 25000
$ws.Range("M2").Value = 24524
$ws.Range("N2").Value = '$/malla 15 kilos'
$ws.Range("O2").Value = 'Hijuelas'
$ws.Range("P2").Value = 1635
$ws.Range("Q2").Value = 15
$ws.Range("R2").Value = 'Hortaliza'

# Row 3
$ws.Range("A3").Value = 9
$ws.Range("B3").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C3").Value = 'Metropolitana'
$ws.Range("D3").Value = 44383
$ws.Range("E3").Value = 13
$ws.Range("F3").Value = 100112035
$ws.Range("G3").Value = 'Bruselas (repollito)'
$ws.Range("H3").Value = 'Sin especificar'
$ws.Range("I3").Value = 'Primera'
$ws.Range("J3").Value = 25
$ws.Range("K3").Value = 13000
$ws.Range("L3").Value = 14000
$ws.Range("M3").Value = 13480
$ws.Range("N3").Value = '$/malla 15 kilos'
$ws.Range("O3").Value = 'Hijuelas'
$ws.Range("P3").Value = 899
$ws.Range("Q3").Value = 15
$ws.Range("R3").Value = 'Hortaliza'

# Row 4
$ws.Range("A4").Value = 9
$ws.Range("B4").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C4").Value = 'Metropolitana'
$ws.Range("D4").Value = 44336
$ws.Range("E4").Value = 13
$ws.Range("F4").Value = 100112035
$ws.Range("G4").Value = 'Bruselas (repollito)'
$ws.Range("H4").Value = 'Sin especificar'
$ws.Range("I4").Value = 'Primera'
$ws.Range("J4").Value = 34
$ws.Range("K4").Value = 24000
$ws.Range("L4").Value = 25000
$ws.Range("M4").Value = 24500
$ws.Range("N4").Value = '$/malla 15 kilos'
$ws.Range("O4").Value = 'Hijuelas'
$ws.Range("P4").Value = 1633
$ws.Range("Q4").Value = 15
$ws.Range("R4").Value = 'Hortaliza'

# Row 5
$ws.Range("A5").Value = 9
$ws.Range("B5").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C5").Value = 'Metropolitana'
$ws.Range("D5").Value = 44390
$ws.Range("E5").Value = 13
$ws.Range("F5").Value = 100112035
$ws.Range("G5").Value = 'Bruselas (repollito)'
$ws.Range("H5").Value = 'Sin especificar'
$ws.Range("I5").Value = 'Primera'
$ws.Range("J5").Value = 34
$ws.Range("K5").Value = 24000
$ws.Range("L5").Value = 25000
$ws.Range("M5").Value = 24500
$ws.Range("N5").Value = '$/malla 15 kilos'
$ws.Range("O5").Value = 'Hijuelas'
$ws.Range("P5").Value = 1633
$ws.Range("Q5").Value = 15
$ws.Range("R5").Value = 'Hortaliza'

# Row 6
$ws.Range("A6").Value = 9
$ws.Range("B6").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C6").Value = 'Metropolitana'
$ws.Range("D6").Value = 44425
$ws.Range("E6").Value = 13
$ws.Range("F6").Value = 100112035
$ws.Range("G6").Value = 'Bruselas (repollito)'
$ws.Range("H6").Value = 'Sin especificar'
$ws.Range("I6").Value = 'Primera'
$ws.Range("J6").Value = 25
$ws.Range("K6").Value = 24000
$ws.Range("L6").Value = 25000
$ws.Range("M6").Value = 24520
$ws.Range("N6").Value = '$/malla 15 kilos'
$ws.Range("O6").Value = 'Hijuelas'
$ws.Range("P6").Value = 1635
$ws.Range("Q6").Value = 15
$ws.Range("R6").Value = 'Hortaliza'

# Row 7
$ws.Range("A7").Value = 9
$ws.Range("B7").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C7").Value = 'Metropolitana'
$ws.Range("D7").Value = 44413
$ws.Range("E7").Value = 13
$ws.Range("F7").Value = 100112035
$ws.Range("G7").Value = 'Bruselas (repollito)'
$ws.Range("H7").Value = 'Sin especificar'
$ws.Range("I7").Value = 'Primera'
$ws.Range("J7").Value = 25
$ws.Range("K7").Value = 24000
$ws.Range("L7").Value = 25000
$ws.Range("M7").Value = 24480
$ws.Range("N7").Value = '$/malla 15 kilos'
$ws.Range("O7").Value = 'Hijuelas'
$ws.Range("P7").Value = 1632
$ws.Range("Q7").Value = 15
$ws.Range("R7").Value = 'Hortaliza'

# Row 8
$ws.Range("A8").Value = 9
$ws.Range("B8").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C8").Value = 'Metropolitana'
$ws.Range("D8").Value = 44418
$ws.Range("E8").Value = 13
$ws.Range("F8").Value = 100112035
$ws.Range("G8").Value = 'Bruselas (repollito)'
$ws.Range("H8").Value = 'Sin especificar'
$ws.Range("I8").Value = 'Primera'
$ws.Range("J8").Value = 16
$ws.Range("K8").Value = 25000
$ws.Range("L8").Value = 26000
$ws.Range("M8").Value = 25500
$ws.Range("N8").Value = '$/malla 15 kilos'
$ws.Range("O8").Value = 'Hijuelas'
$ws.Range("P8").Value = 1700
$ws.Range("Q8").Value = 15
$ws.Range("R8").Value = 'Hortaliza'

# Row 9
$ws.Range("A9").Value = 9
$ws.Range("B9").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C9").Value = 'Metropolitana'
$ws.Range("D9").Value = 44343
$ws.Range("E9").Value = 13
$ws.Range("F9").Value = 100112035
$ws.Range("G9").Value = 'Bruselas (repollito)'
$ws.Range("H9").Value = 'Sin especificar'
$ws.Range("I9").Value = 'Primera'
$ws.Range("J9").Value = 26
$ws.Range("K9").Value = 23000
$ws.Range("L9").Value = 24000
$ws.Range("M9").Value = 23500
$ws.Range("N9").Value = '$/malla 15 kilos'
$ws.Range("O9").Value = 'Hijuelas'
$ws.Range("P9").Value = 1567
$ws.Range("Q9").Value = 15
$ws.Range("R9").Value = 'Hortaliza'

# Row 10
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C10").Value = 'Metropolitana'
$ws.Range("D10").Value = 44406
$ws.Range("E10").Value = 13
$ws.Range("F10").Value = 100112035
$ws.Range("G10").Value = 'Bruselas (repollito)'
$ws.Range("H10").Value = 'Sin especificar'
$ws.Range("I10").Value = 'Primera'
$ws.Range("J10").Value = 25
$ws.Range("K10").Value = 24000
$ws.Range("L10").Value = 25000
$ws.Range("M10").Value = 24520
$ws.Range("N10").Value = '$/malla 15 kilos'
$ws.Range("O10").Value = 'Hijuelas'
$ws.Range("P10").Value = 1635
$ws.Range("Q10").Value = 15
$ws.Range("R10").Value = 'Hortaliza'

# Row 11
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C11").Value = 'Metropolitana'
$ws.Range("D11").Value = 44385
$ws.Range("E11").Value = 13
$ws.Range("F11").Value = 100112035
$ws.Range("G11").Value = 'Bruselas (repollito)'
$ws.Range("H11").Value = 'Sin especificar'
$ws.Range("I11").Value = 'Primera'
$ws.Range("J11").Value = 25
$ws.Range("K11").Value = 14000
$ws.Range("L11").Value = 15000
$ws.Range("M11").Value = 14480
$ws.Range("N11").Value = '$/malla 15 kilos'
$ws.Range("O11").Value = 'Hijuelas'
$ws.Range("P11").Value = 965
$ws.Range("Q11").Value = 15
$ws.Range("R11").Value = 'Hortaliza'

# Row 12
$ws.Range("A12").Value = 9
$ws.Range("B12").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C12").Value = 'Metropolitana'
$ws.Range("D12").Value = 44385
$ws.Range("E12").Value = 13
$ws.Range("F12").Value = 100112035
$ws.Range("G12").Value = 'Bruselas (repollito)'
$ws.Range("H12").Value = 'Sin especificar'
$ws.Range("I12").Value = 'Segunda'
$ws.Range("J12").Value = 16
$ws.Range("K12").Value = 12000
$ws.Range("L12").Value = 12000
$ws.Range("M12").Value = 12000
$ws.Range("N12").Value = '$/malla 15 kilos'
$ws.Range("O12").Value = 'Hijuelas'
$ws.Range("P12").Value = 800
$ws.Range("Q12").Value = 15
$ws.Range("R12").Value = 'Hortaliza'

# Row 13
$ws.Range("A13").Value = 9
$ws.Range("B13").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C13").Value = 'Metropolitana'
$ws.Range("D13").Value = 44397
$ws.Range("E13").Value = 13
$ws.Range("F13").Value = 100112035
$ws.Range("G13").Value = 'Bruselas (repollito)'
$ws.Range("H13").Value = 'Sin especificar'
$ws.Range("I13").Value = 'Primera'
$ws.Range("J13").Value = 34
$ws.Range("K13").Value = 23000
$ws.Range("L13").Value = 24000
$ws.Range("M13").Value = 23500
$ws.Range("N13").Value = '$/malla 15 kilos'
$ws.Range("O13").Value = 'Hijuelas'
$ws.Range("P13").Value = 1567
$ws.Range("Q13").Value = 15
$ws.Range("R13").Value = 'Hortaliza'

# Row 14
$ws.Range("A14").Value = 9
$ws.Range("B14").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C14").Value = 'Metropolitana'
$ws.Range("D14").Value = 44421
$ws.Range("E14").Value = 13
$ws.Range("F14").Value = 100112035
$ws.Range("G14").Value = 'Bruselas (repollito)'
$ws.Range("H14").Value = 'Sin especificar'
$ws.Range("I14").Value = 'Primera'
$ws.Range("J14").Value = 18
$ws.Range("K14").Value = 24000
$ws.Range("L14").Value = 25000
$ws.Range("M14").Value = 24500
$ws.Range("N14").Value = '$/malla 15 kilos'
$ws.Range("O14").Value = 'Hijuelas'
$ws.Range("P14").Value = 1633
$ws.Range("Q14").Value = 15
$ws.Range("R14").Value = 'Hortaliza'

# Row 15
$ws.Range("A15").Value = 9
$ws.Range("B15").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C15").Value = 'Metropolitana'
$ws.Range("D15").Value = 44341
$ws.Range("E15").Value = 13
$ws.Range("F15").Value = 100112035
$ws.Range("G15").Value = 'Bruselas (repollito)'
$ws.Range("H15").Value = 'Sin especificar'
$ws.Range("I15").Value = 'Primera'
$ws.Range("J15").Value = 36
$ws.Range("K15").Value = 24000
$ws.Range("L15").Value = 25000
$ws.Range("M15").Value = 24500
$ws.Range("N15").Value = '$/malla 15 kilos'
$ws.Range("O15").Value = 'Hijuelas'
$ws.Range("P15").Value = 1633
$ws.Range("Q15").Value = 15
$ws.Range("R15").Value = 'Hortaliza'

# Row 16
$ws.Range("A16").Value = 9
$ws.Range("B16").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C16").Value = 'Metropolitana'
$ws.Range("D16").Value = 44432
$ws.Range("E16").Value = 13
$ws.Range("F16").Value = 100112035
$ws.Range("G16").Value = 'Bruselas (repollito)'
$ws.Range("H16").Value = 'Sin especificar'
$ws.Range("I16").Value = 'Primera'
$ws.Range("J16").Value = 34
$ws.Range("K16").Value = 24000
$ws.Range("L16").Value = 25000
$ws.Range("M16").Value = 24500
$ws.Range("N16").Value = '$/malla 15 kilos'
$ws.Range("O16").Value = 'Hijuelas'
$ws.Range("P16").Value = 1633
$ws.Range("Q16").Value = 15
$ws.Range("R16").Value = 'Hortaliza'

# Row 17
$ws.Range("A17").Value = 9
$ws.Range("B17").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C17").Value = 'Metropolitana'
$ws.Range("D17").Value = 44351
$ws.Range("E17").Value = 13
$ws.Range("F17").Value = 100112035
$ws.Range("G17").Value = 'Bruselas (repollito)'
$ws.Range("H17").Value = 'Sin especificar'
$ws.Range("I17").Value = 'Primera'
$ws.Range("J17").Value = 34
$ws.Range("K17").Value = 24000
$ws.Range("L17").Value = 25000
$ws.Range("M17").Value = 24500
$ws.Range("N17").Value = '$/malla 15 kilos'
$ws.Range("O17").Value = 'Hijuelas'
$ws.Range("P17").Value = 1633
$ws.Range("Q17").Value = 15
$ws.Range("R17").Value = 'Hortaliza'

# Row 18
$ws.Range("A18").Value = 9
$ws.Range("B18").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C18").Value = 'Metropolitana'
$ws.Range("D18").Value = 44446
$ws.Range("E18").Value = 13
$ws.Range("F18").Value = 100112035
$ws.Range("G18").Value = 'Bruselas (repollito)'
$ws.Range("H18").Value = 'Sin especificar'
$ws.Range("I18").Value = 'Primera'
$ws.Range("J18").Value = 34
$ws.Range("K18").Value = 24000
$ws.Range("L18").Value = 25000
$ws.Range("M18").Value = 24500
$ws.Range("N18").Value = '$/malla 15 kilos'
$ws.Range("O18").Value = 'Hijuelas'
$ws.Range("P18").Value = 1633
$ws.Range("Q18").Value = 15
$ws.Range("R18").Value = 'Hortaliza'

# Row 19
$ws.Range("A19").Value = 9
$ws.Range("B19").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C19").Value = 'Metropolitana'
$ws.Range("D19").Value = 44411
$ws.Range("E19").Value = 13
$ws.Range("F19").Value = 100112035
$ws.Range("G19").Value = 'Bruselas (repollito)'
$ws.Range("H19").Value = 'Sin especificar'
$ws.Range("I19").Value = 'Primera'
$ws.Range("J19").Value = 34
$ws.Range("K19").Value = 25000
$ws.Range("L19").Value = 26000
$ws.Range("M19").Value = 25500
$ws.Range("N19").Value = '$/malla 15 kilos'
$ws.Range("O19").Value = 'Hijuelas'
$ws.Range("P19").Value = 1700
$ws.Range("Q19").Value = 15
$ws.Range("R19").Value = 'Hortaliza'

# Row 20
$ws.Range("A20").Value = 9
$ws.Range("B20").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C20").Value = 'Metropolitana'
$ws.Range("D20").Value = 44329
$ws.Range("E20").Value = 13
$ws.Range("F20").Value = 100112035
$ws.Range("G20").Value = 'Bruselas (repollito)'
$ws.Range("H20").Value = 'Sin especificar'
$ws.Range("I20").Value = 'Primera'
$ws.Range("J20").Value = 25
$ws.Range("K20").Value = 23000
$ws.Range("L20").Value = 23000
$ws.Range("M20").Value = 23000
$ws.Range("N20").Value = '$/malla 15 kilos'
$ws.Range("O20").Value = 'Hijuelas'
$ws.Range("P20").Value = 1533
$ws.Range("Q20").Value = 15
$ws.Range("R20").Value = 'Hortaliza'

# Row 21
$ws.Range("A21").Value = 9
$ws.Range("B21").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C21").Value = 'Metropolitana'
$ws.Range("D21").Value = 44428
$ws.Range("E21").Value = 13
$ws.Range("F21").Value = 100112035
$ws.Range("G21").Value = 'Bruselas (repollito)'
$ws.Range("H21").Value = 'Sin especificar'
$ws.Range("I21").Value = 'Primera'
$ws.Range("J21").Value = 16
$ws.Range("K21").Value = 25000
$ws.Range("L21").Value = 26000
$ws.Range("M21").Value = 25500
$ws.Range("N21").Value = '$/malla 15 kilos'
$ws.Range("O21").Value = 'Hijuelas'
$ws.Range("P21").Value = 1700
$ws.Range("Q21").Value = 15
$ws.Range("R21").Value = 'Hortaliza'

# Row 22
$ws.Range("A22").Value = 9
$ws.Range("B22").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C22").Value = 'Metropolitana'
$ws.Range("D22").Value = 44442
$ws.Range("E22").Value = 13
$ws.Range("F22").Value = 100112035
$ws.Range("G22").Value = 'Bruselas (repollito)'
$ws.Range("H22").Value = 'Sin especificar'
$ws.Range("I22").Value = 'Primera'
$ws.Range("J22").Value = 28
$ws.Range("K22").Value = 24000
$ws.Range("L22").Value = 25000
$ws.Range("M22").Value = 24500
$ws.Range("N22").Value = '$/malla 15 kilos'
$ws.Range("O22").Value = 'Hijuelas'
$ws.Range("P22").Value = 1633
$ws.Range("Q22").Value = 15
$ws.Range("R22").Value = 'Hortaliza'

# Row 23
$ws.Range("A23").Value = 9
$ws.Range("B23").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C23").Value = 'Metropolitana'
$ws.Range("D23").Value = 44435
$ws.Range("E23").Value = 13
$ws.Range("F23").Value = 100112035
$ws.Range("G23").Value = 'Bruselas (repollito)'
$ws.Range("H23").Value = 'Sin especificar'
$ws.Range("I23").Value = 'Primera'
$ws.Range("J23").Value = 34
$ws.Range("K23").Value = 24000
$ws.Range("L23").Value = 25000
$ws.Range("M23").Value = 24500
$ws.Range("N23").Value = '$/malla 15 kilos'
$ws.Range("O23").Value = 'Hijuelas'
$ws.Range("P23").Value = 1633
$ws.Range("Q23").Value = 15
$ws.Range("R23").Value = 'Hortaliza'
$ws.Range("D23").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 24
$ws.Range("A24").Value = 9
$ws.Range("B24").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C24").Value = 'Metropolitana'
$ws.Range("D24").Value = 44400
$ws.Range("E24").Value = 13
$ws.Range("F24").Value = 100112035
$ws.Range("G24").Value = 'Bruselas (repollito)'
$ws.Range("H24").Value = 'Sin especificar'
$ws.Range("I24").Value = 'Primera'
$ws.Range("J24").Value = 16
$ws.Range("K24").Value = 24000
$ws.Range("L24").Value = 25000
$ws.Range("M24").Value = 24500
$ws.Range("N24").Value = '$/malla 15 kilos'
$ws.Range("O24").Value = 'Hijuelas'
$ws.Range("P24").Value = 1633
$ws.Range("Q24").Value = 15
$ws.Range("R24").Value = 'Hortaliza'
$ws.Range("D24").NumberFormat = "YYYY-MM-DD HH:MM:SS"
